$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("A2:D18")
$key1 = $ws.Range("A2:A18")
$rng.Sort($key1, 1)
